$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the end of the first paragraph.
# In the edited document it moves to the end of the new last paragraph, so
# remove it here; it will be re-created (by id/name) in the inserted markup.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Append a blank paragraph followed by a new paragraph reading
# "Modificacion de la portatil" (with the spell-check proofErr markers
# around "Modificacion", as Word would leave them) and carrying the
# relocated _GoBack bookmark at its end.
$r = $d.Content
$r.Collapse(0)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body>' +
                    '<w:p/>' +
                    '<w:p>' +
                        '<w:proofErr w:type="spellStart"/>' +
                        '<w:r><w:t>Modificacion</w:t></w:r>' +
                        '<w:proofErr w:type="spellEnd"/>' +
                        '<w:r><w:t xml:space="preserve"> de la portatil</w:t></w:r>' +
                        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
                        '<w:bookmarkEnd w:id="0"/>' +
                    '</w:p>' +
                '</w:body>' +
            '</w:document>' +
        '</pkg:xmlData>' +
    '</pkg:part>' +
'</pkg:package>'

$r.InsertXML($xml)
